# Generate Report for handoff
# Updates the localization-status workbook: the source markdown file was
# regenerated under a new GUID, and the zh-cn / de-de handoff attempts are
# now reported as failed/ignored instead of successfully handed off.

$wb = $excel.ActiveWorkbook

$oldFile = "b4327c94-b673-4614-b729-626121798d22.md"
$newFile = "a2f13a79-99af-45fb-b9b8-693c2a92bee9.md"

$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/476b3a78e894ccd09caf4290fffd22cdf4722ab2/e2e/" + $newFile
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/476b3a78e894ccd09caf4290fffd22cdf4722ab2/.localization-config"

# ---------------------------------------------------------------------
# Overview sheet: just the file name + status text changes.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Per-locale sheets: the handoff attempt failed, so the "Latest Handoff
# File" link is cleared, the handoff datetime resets to the sentinel
# value, and the reason flips from Include to Ignored.
# ---------------------------------------------------------------------
function Update-LocaleSheet($sheetName, $resetDate) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value = $newFile
    $ws.Range("B2").Value = $newStatus
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = $resetDate
    $ws.Range("G2").Value = $resetDate
    $ws.Range("H2").Value = "Ignored"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $newFile)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config")
}

Update-LocaleSheet "zh-cn" "0001-01-01 00:00:00"
Update-LocaleSheet "de-de" "0001-01-01 00:00:00"
